$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 71, shifting existing rows 71-119 down to 72-120.
$ws.Rows(71).Insert()

# Populate the newly inserted row 71 with the new record.
$ws.Cells.Item(71, 1).Value = 8
$ws.Cells.Item(71, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(71, 3).Value = "Coquimbo"
$ws.Cells.Item(71, 4).Value = 44977
$ws.Cells.Item(71, 5).Value = 4
$ws.Cells.Item(71, 6).Value = 100112030
$ws.Cells.Item(71, 7).Value = "Poroto granado"
$ws.Cells.Item(71, 8).Value = "Sin especificar"
$ws.Cells.Item(71, 9).Value = "Primera"
$ws.Cells.Item(71, 10).Value = 400
$ws.Cells.Item(71, 11).Value = 37000
$ws.Cells.Item(71, 12).Value = 38000
$ws.Cells.Item(71, 13).Value = 37500
$ws.Cells.Item(71, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(71, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(71, 16).Value = 1500
$ws.Cells.Item(71, 17).Value = 25
$ws.Cells.Item(71, 18).Value = "Hortaliza"
